$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K")

# Row 4 - Inventory
$ws.Range("B4").Value = 1284000000.0
$ws.Range("C4").Value = 1263000000.0
$ws.Range("D4").Value = 1230000000.0
$ws.Range("E4").Value = 1189000000.0
$ws.Range("F4").Value = 1226000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 2471000000.0
$ws.Range("C14").Value = 2449000000.0
$ws.Range("D14").Value = 2393000000.0
$ws.Range("E14").Value = 2329000000.0
$ws.Range("F14").Value = 2387000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = 308000000.0
$ws.Range("C22").Value = 368000000.0
$ws.Range("D22").Value = 374000000.0
$ws.Range("E22").Value = 370000000.0
$ws.Range("F22").Value = 365000000.0
